$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 259
$ws.Range("C3").Value = 157811
$ws.Range("C4").Value = 148873
$ws.Range("C7").Value = 5.66
$ws.Range("C8").Value = 63.87
